$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")
Write-Host $ws.Range("A1").Value
